# Fruta / hortaliza, semanal
# Insert 3 new weekly data rows at the top of the existing "Frutilla" data
# block (rows 869-871), pushing the rest of the data down by 3 rows. The
# previously-last 3 rows end up duplicated at the very end of the sheet as
# a consequence of the shift (rows 909-911), matching the new dimension
# A1:T911.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 869:908 down to 872:911, inserting 3 blank rows at 869.
$ws.Rows("869:871").Insert()

# New row 869: Especial, Provincia de Melipilla, semana del 2022-05-09 (44706)
$ws.Range("A869").Value2 = 9
$ws.Range("B869").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C869").Value = "Metropolitana"
$ws.Range("D869").Value2 = 44706
$ws.Range("E869").Value2 = 13
$ws.Range("F869").Value = "Fruta"
$ws.Range("G869").Value2 = 100101
$ws.Range("H869").Value = "Berries"
$ws.Range("I869").Value2 = 100112025
$ws.Range("J869").Value = "Frutilla"
$ws.Range("K869").Value = "Sin especificar"
$ws.Range("L869").Value = "Especial"
$ws.Range("M869").Value2 = 330
$ws.Range("N869").Value2 = 12000
$ws.Range("O869").Value2 = 12000
$ws.Range("P869").Value2 = 12000
$ws.Range("Q869").Value = "$/bandeja 7 kilos"
$ws.Range("R869").Value = "Provincia de Melipilla"
$ws.Range("S869").Value2 = 1714
$ws.Range("T869").Value2 = 7

# New row 870: Primera, Provincia de Melipilla
$ws.Range("A870").Value2 = 9
$ws.Range("B870").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C870").Value = "Metropolitana"
$ws.Range("D870").Value2 = 44706
$ws.Range("E870").Value2 = 13
$ws.Range("F870").Value = "Fruta"
$ws.Range("G870").Value2 = 100101
$ws.Range("H870").Value = "Berries"
$ws.Range("I870").Value2 = 100112025
$ws.Range("J870").Value = "Frutilla"
$ws.Range("K870").Value = "Sin especificar"
$ws.Range("L870").Value = "Primera"
$ws.Range("M870").Value2 = 310
$ws.Range("N870").Value2 = 10000
$ws.Range("O870").Value2 = 10000
$ws.Range("P870").Value2 = 10000
$ws.Range("Q870").Value = "$/bandeja 7 kilos"
$ws.Range("R870").Value = "Provincia de Melipilla"
$ws.Range("S870").Value2 = 1429
$ws.Range("T870").Value2 = 7

# New row 871: Segunda, Provincia de Melipilla
$ws.Range("A871").Value2 = 9
$ws.Range("B871").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C871").Value = "Metropolitana"
$ws.Range("D871").Value2 = 44706
$ws.Range("E871").Value2 = 13
$ws.Range("F871").Value = "Fruta"
$ws.Range("G871").Value2 = 100101
$ws.Range("H871").Value = "Berries"
$ws.Range("I871").Value2 = 100112025
$ws.Range("J871").Value = "Frutilla"
$ws.Range("K871").Value = "Sin especificar"
$ws.Range("L871").Value = "Segunda"
$ws.Range("M871").Value2 = 280
$ws.Range("N871").Value2 = 7000
$ws.Range("O871").Value2 = 7000
$ws.Range("P871").Value2 = 7000
$ws.Range("Q871").Value = "$/bandeja 7 kilos"
$ws.Range("R871").Value = "Provincia de Melipilla"
$ws.Range("S871").Value2 = 1000
$ws.Range("T871").Value2 = 7
